$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.214.77'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '3.373.28'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '406.39'
$ws.Range('E5').Value = '  -1.52%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.16'
$ws.Range('E6').Value = '  +9.42%  '
$ws.Range('E7').Value = '  +1.56%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.670'
$ws.Range('E9').Value = '  +4.39%  '
$ws.Range('E10').Value = '  -3.70%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '42.95'
$ws.Range('E11').Value = '  +3.97%  '
$ws.Range('E12').Value = '  -0.96%  '
$ws.Range('D13').Value = '3.888.88'
$ws.Range('E13').Value = '  -1.40%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '8.36'
$ws.Range('E14').Value = '  -0.96%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '19.59'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('D16').Value = '3.372.95'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').Value = '61.164.80'
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.01'
$ws.Range('E19').Value = '  +1.33%  '
$ws.Range('E20').Value = '  +3.09%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '3.20'
$ws.Range('E21').Value = '  -3.26%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '83.48'
$ws.Range('E22').Value = '  +7.72%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '312.75'
$ws.Range('E23').Value = '  +4.15%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.78'
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.78'
$ws.Range('E26').Value = '  +11.64%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.32'
$ws.Range('E27').Value = '  +9.10%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '29.42'
$ws.Range('E28').Value = '  -3.75%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.52'
$ws.Range('E29').Value = '  -6.37%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '11.30'
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('B33').Value = 'Dai'
$ws.Range('C33').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '41.48'
$ws.Range('E34').Value = '  -1.11%  '
$ws.Range('E35').Value = '  -2.13%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0480'
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '52.09'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.997'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('E39').Value = '  -2.75%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.93'
$ws.Range('E40').Value = '  -2.27%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '137.77'
$ws.Range('E41').Value = '  +2.60%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.97'
$ws.Range('E42').Value = '  +0.69%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.123'
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.293'
$ws.Range('E44').Value = '  +4.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.03'
$ws.Range('E45').Value = '  +3.44%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '16.66'
$ws.Range('E46').Value = '  -4.48%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.22'
$ws.Range('E47').Value = '  +1.94%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '21.39'
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('D49').Value = '2.123.15'
$ws.Range('E49').Value = '  -3.45%  '
$ws.Range('E50').Value = '  -4.90%  '
$ws.Range('E51').Value = '  -0.18%  '
